# Update data through 2022-04-07: rename the "through March 29" column/sheet
# to "through March 30" and bump the affected neighborhood counts for
# March 2022 (column B), plus a handful of other corrected cells across
# the historical columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab and header label (shared string) to reflect the new "through" date.
$ws.Name = "Through 2022-03-30"
$ws.Range("B1").Value = "March 2022 (through March 30)"

# Row 3 - Austin
$ws.Range("N3").Value = 3
$ws.Range("Q3").Value = 3

# Row 5 - Garfield Park
$ws.Range("B5").Value = 8
$ws.Range("E5").Value = 5

# Row 9 - Chicago Lawn
$ws.Range("B9").Value = 4

# Row 11 - Englewood
$ws.Range("H11").Value = 4

# Row 14 - West Town (new value)
$ws.Range("Q14").Value = 1

# Row 15 - Humboldt Park
$ws.Range("E15").Value = 2

# Row 16 - Little Italy, UIC
$ws.Range("K16").Value = 2

# Row 17 - Auburn Gresham
$ws.Range("B17").Value = 2

# Row 23 - Wicker Park
$ws.Range("B23").Value = 3

# Row 43 - Hermosa (new value)
$ws.Range("Q43").Value = 1

# Row 44 - Grand Boulevard (new values)
$ws.Range("B44").Value = 1
$ws.Range("W44").Value = 1

# Row 47 - Brighton Park (new value)
$ws.Range("N47").Value = 1

# Row 85 - South Chicago
$ws.Range("E85").Value = 2
$ws.Range("H85").Value = 2
